$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49; this shifts the existing rows 49-139
# down to 50-140 (carrying all of their values/formatting with them).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R keep the same constant /
# repeated values used throughout the table; D, J, K, L, M, P are the
# new data points for this entry.
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44868
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 100114002
$ws.Range("G49").Value = "Camote"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 170
$ws.Range("K49").Value = 25000
$ws.Range("L49").Value = 26000
$ws.Range("M49").Value = 25529
$ws.Range("N49").Value = "$/malla 20 kilos"
$ws.Range("O49").Value = "Perú"
$ws.Range("P49").Value = 1276
$ws.Range("Q49").Value = 20
$ws.Range("R49").Value = "Hortaliza"
